# Undo Jason's overwrite of class materials starter code/slides.
#
# Jason Yoder had added a highlighted "Today's Attendance password /
# merging" textbox to the title slide. This restores the slide to the
# prior (Cameron Dorsey) version by removing that textbox, and restores
# the cached "today" date captions on the slide master / layouts (and
# handout/notes masters) back to their prior cached values.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Remove the textbox Jason added to the title slide (slide id 304 ==
#    the first slide in the deck): shape id 3, name "TextBox 2",
#    containing the highlighted "Today's Attendance password / merging"
#    text.
# ---------------------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
for ($i = $titleSlide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $titleSlide.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 2" -and $shp.Id -eq 3) {
        $shp.Delete()
    }
}

# ---------------------------------------------------------------------
# 2. Restore the cached date-placeholder text (PowerPoint caches the
#    evaluated text of the "Update automatically" date field; it had
#    drifted to the 9/3/2023 save date and needs restoring to the prior
#    11/6/2022 save date) across the slide master, every slide layout,
#    the handout master, and the notes master.
# ---------------------------------------------------------------------
function Set-DatePlaceholderText {
    param($shapes, [string]$text)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $text
        }
    }
}

$longDate = "Sunday, November 6, 2022"
$shortDate = "11/6/22"

# Slide master
Set-DatePlaceholderText $p.SlideMaster.Shapes $longDate

# Every slide layout off the (single) slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes $longDate
}

# Handout master + notes master use the short m/d/yy cached format
Set-DatePlaceholderText $p.HandoutMaster.Shapes $shortDate
Set-DatePlaceholderText $p.NotesMaster.Shapes $shortDate
